# Weekly price-data update: insert a new week's worth of "Frutilla" (strawberry)
# price rows (Especial / Primera / Segunda / Tercera) at the top of the detail
# rows (rows 584-587), pushing all the existing rows down by 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at position 584 (shifts old rows 584..682 down to 588..686)
$ws.Rows("584:587").Insert()

# Common (unchanging) column values for this market/product block
$mercadoId   = 6
$mercado     = "Mercado Mayorista Lo Valledor de Santiago"
$region      = "Metropolitana"
$fecha       = 44505
$codreg      = 13
$tipo        = "Fruta"
$productoId  = 100101
$producto    = "Berries"
$categoriaId = 100112025
$categoria   = "Frutilla"
$variedad    = "Sin especificar"
$unidad      = "`$/bandeja 7 kilos"
$origen      = "Provincia de Melipilla"
$kgUnidad    = 7

function Set-FrutillaRow($Row, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg) {
    $ws.Cells.Item($Row, 1).Value  = $mercadoId
    $ws.Cells.Item($Row, 2).Value  = $mercado
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $fecha
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $tipo
    $ws.Cells.Item($Row, 7).Value  = $productoId
    $ws.Cells.Item($Row, 8).Value  = $producto
    $ws.Cells.Item($Row, 9).Value  = $categoriaId
    $ws.Cells.Item($Row, 10).Value = $categoria
    $ws.Cells.Item($Row, 11).Value = $variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $unidad
    $ws.Cells.Item($Row, 18).Value = $origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $kgUnidad
}

Set-FrutillaRow 584 "Especial" 2950 6000 7000 6500 929
Set-FrutillaRow 585 "Primera"  2800 4500 5500 4911 702
Set-FrutillaRow 586 "Segunda"  900  3000 4000 3500 500
Set-FrutillaRow 587 "Tercera"  500  2000 2500 2250 321
